$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331, shifting existing rows 331-417 down to 332-418.
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row 331 with the new weekly data point.
$ws.Cells.Item(331, 1).Value = 3
$ws.Cells.Item(331, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 44841
$ws.Cells.Item(331, 5).Value = 5
$ws.Cells.Item(331, 6).Value = 100112012
$ws.Cells.Item(331, 7).Value = "Espinaca"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 120
$ws.Cells.Item(331, 11).Value = 4000
$ws.Cells.Item(331, 12).Value = 4000
$ws.Cells.Item(331, 13).Value = 4000
$ws.Cells.Item(331, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(331, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(331, 16).Value = 1333
$ws.Cells.Item(331, 17).Value = 3
$ws.Cells.Item(331, 18).Value = "Hortaliza"
